# "Add cantrals by cantons"
# Restructure Sheet1: move from a two-row header (with a units row) and a
# data block that has some "orphan" rows, into a single-row header
# (idx, idx2, Name, Date Start, Date End, (m3/s), (MW1), (MW2),
#  (GWh) Winter, (GWh) Summer, (GWh) Year) followed directly by 6 data
# rows (one per power plant / "central").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Start from a clean sheet - clear everything currently used.
$ws.Cells.Clear()

# ---- Header row -------------------------------------------------------
$headers = @("idx", "idx2", "Name", "Date Start", "Date End", "(m3/s)", "(MW1)", "(MW2)", "(GWh) Winter", "(GWh) Summer", "(GWh) Year")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# The F1:K1 header cells carry the "units" look (font id 1, Arial 9,
# default/general number format).
$unitsHeaderRange = $ws.Range($ws.Cells.Item(1, 6), $ws.Cells.Item(1, 11))
$unitsHeaderRange.Font.Name = "Arial"
$unitsHeaderRange.Font.Size = 9

# ---- Data rows ----------------------------------------------------------
# idx, idx2, Name, Date Start, Date End, (m3/s), (MW1), (MW2), (GWh) Winter, (GWh) Summer, (GWh) Year
$data = @(
    @(1, 108900, "Wunderklingen", 1895, 1968, 5.5,   0.42,  0.41,               1.4,   1,     2.4),
    @(2, 106300, "Engeweiher",    1909, 1993, 4,     5,     5,                  $null, $null, $null),
    @(3, 108700, "Eglisau",       1920, 2012, 500,   14.91, 16.920000000000002, 47.38, 53.74, 101.12),
    @(4, 106400, "Neuhausen",     1951, 2011, 29.9,  2.8,   2.4500000000000002, 10.35, 10.4,  20.75),
    @(5, 106500, "Rheinau",       1956, 2005, 400,   2.98,  2.92,               6.39,  13.19, 19.579999999999998),
    @(6, 106200, "Schaffhausen",  1964, $null, 500,  22.57, 19.84,              62.06, 73.64, 135.69999999999999)
)

# Column layout:
#   A,B,D,E -> integer id / year values
#   C       -> plant name (text)
#   F:K     -> 2 decimal measurement values
$intCols = @(1, 2, 4, 5)
$decimalCols = @(6, 7, 8, 9, 10, 11)

$rowIndex = 2
foreach ($row in $data) {
    for ($col = 1; $col -le $row.Length; $col++) {
        $value = $row[$col - 1]
        if ($null -eq $value) {
            continue
        }
        $cell = $ws.Cells.Item($rowIndex, $col)
        $cell.Value = $value
        $cell.Font.Name = "Arial"
        $cell.Font.Size = 9
        if ($intCols -contains $col) {
            $cell.NumberFormat = "0"
        } elseif ($decimalCols -contains $col) {
            $cell.NumberFormat = "0.00"
        }
    }
    $rowIndex++
}

# ---- Selection / view ----------------------------------------------------
$ws.Range("A2:K2").Select()
